$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting current rows 18-19 down to 19-20.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44943
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100114007
$ws.Range("G18").Value = "Jengibre"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 350
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14429
$ws.Range("N18").Value = "$/caja 13 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 1110
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "Hortaliza"
